# Auto-generated Excel COM-interop script to apply the OOXML diff described in the commit.
# For each affected sheet/row, numeric cells in columns H:N (computed price/profit columns)
# are updated to match the target values. Cells that were removed entirely in the diff are
# cleared with ClearContents() so they no longer serialize as <c> elements; cells that were
# newly introduced in the diff are written for the first time.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 179.6
$ws.Range("I2").Value = 177.44444
$ws.Range("J2").Value = 199
$ws.Range("K2").Value = 177.44444
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = -64.44443999999999
$ws.Range("N2").Value = -425
# Row 12
$ws.Range("H12").Value = 419.6
$ws.Range("I12").Value = 199.66667
$ws.Range("K12").Value = 199.66667
$ws.Range("M12").Value = -29.66667000000001
# Row 74
$ws.Range("H74").Value = 10319.667
$ws.Range("J74").Value = 11000
$ws.Range("L74").Value = 11000
$ws.Range("N74").Value = -12872
# Row 77
$ws.Range("H77").Value = 10319.667
$ws.Range("J77").Value = 11000
$ws.Range("L77").Value = 55000
$ws.Range("N77").Value = -64360
# Row 80
$ws.Range("H80").Value = 2917.8845
$ws.Range("I80").Value = 1480.2222
$ws.Range("J80").Value = 3679
$ws.Range("K80").Value = 4440.6666
$ws.Range("L80").Value = 11037
$ws.Range("M80").Value = -3442.6666
$ws.Range("N80").Value = -13033
# Row 83
$ws.Range("H83").Value = 2917.8845
$ws.Range("I83").Value = 1480.2222
$ws.Range("J83").Value = 3679
$ws.Range("K83").Value = 13321.9998
$ws.Range("L83").Value = 33111
$ws.Range("M83").Value = -8329.9998
$ws.Range("N83").Value = -43095
# Row 98
$ws.Range("H98").Value = 2120.125
$ws.Range("I98").Value = 2888.75
$ws.Range("K98").Value = 2888.75
$ws.Range("M98").Value = -1390.75
# Row 100
$ws.Range("H100").Value = 5845.48
$ws.Range("I100").Value = 4650.364
$ws.Range("J100").Value = 6784.5
$ws.Range("K100").Value = 4650.364
$ws.Range("L100").Value = 6784.5
$ws.Range("M100").Value = -4109.364
$ws.Range("N100").Value = -7866.5
# Row 101
$ws.Range("H101").Value = 570
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 122
$ws.Range("H122").Value = 2120.125
$ws.Range("I122").Value = 2888.75
$ws.Range("K122").Value = 8666.25
$ws.Range("M122").Value = -6216.25
# Row 137
$ws.Range("H137").Value = 2139.6223
$ws.Range("I137").Value = 1985.3143
$ws.Range("K137").Value = 5955.9429
$ws.Range("M137").Value = -3405.9429

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 282.66666
$ws.Range("I4").Value = 282.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 282.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -166.66666
$ws.Range("N4").ClearContents()
# Row 45
$ws.Range("H45").Value = 3729.524
$ws.Range("I45").Value = 2475.5833
$ws.Range("J45").Value = 5401.4443
$ws.Range("K45").Value = 2475.5833
$ws.Range("L45").Value = 5401.4443
$ws.Range("M45").Value = -2098.5833
$ws.Range("N45").Value = -6155.4443
# Row 61
$ws.Range("H61").Value = 1824.7778
$ws.Range("I61").Value = 1824.7778
$ws.Range("K61").Value = 1824.7778
$ws.Range("M61").Value = -1612.7778
# Row 110
$ws.Range("H110").Value = 6636.273
$ws.Range("I110").Value = 5799.9
$ws.Range("K110").Value = 5799.9
$ws.Range("M110").Value = -3754.9
# Row 132
$ws.Range("H132").Value = 2508.8462
$ws.Range("I132").Value = 2197.0833
$ws.Range("K132").Value = 6591.249899999999
$ws.Range("M132").Value = -4061.249899999999
# Row 136
$ws.Range("H136").Value = 1824.7778
$ws.Range("I136").Value = 1824.7778
$ws.Range("K136").Value = 5474.3334
$ws.Range("M136").Value = -2924.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2733.9285
$ws.Range("I20").Value = 2572
$ws.Range("J20").Value = 2798.7
$ws.Range("K20").Value = 2572
$ws.Range("L20").Value = 2798.7
$ws.Range("M20").Value = -2325
$ws.Range("N20").Value = -3292.7
# Row 93
$ws.Range("H93").Value = 32500
$ws.Range("J93").Value = 32500
$ws.Range("L93").Value = 32500
$ws.Range("N93").Value = -36244
# Row 96
$ws.Range("H96").Value = 12500
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
# Row 103
$ws.Range("H103").Value = 10999.333
$ws.Range("J103").Value = 10999.333
$ws.Range("L103").Value = 10999.333
$ws.Range("N103").Value = -13343.333
# Row 105
$ws.Range("H105").Value = 2074.625
$ws.Range("I105").Value = 1839.4
$ws.Range("J105").Value = 2466.6667
$ws.Range("K105").Value = 1839.4
$ws.Range("L105").Value = 2466.6667
$ws.Range("M105").Value = -92.40000000000009
$ws.Range("N105").Value = -5960.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1991.5
$ws.Range("I31").Value = 2051.077
$ws.Range("K31").Value = 2051.077
$ws.Range("M31").Value = -1756.077
# Row 34
$ws.Range("H34").Value = 1991.5
$ws.Range("I34").Value = 2051.077
$ws.Range("K34").Value = 2051.077
$ws.Range("M34").Value = -1849.077
# Row 122
$ws.Range("H122").Value = 5049.95
$ws.Range("I122").Value = 4463.125
$ws.Range("J122").Value = 5441.1665
$ws.Range("K122").Value = 13389.375
$ws.Range("L122").Value = 16323.4995
$ws.Range("M122").Value = -10939.375
$ws.Range("N122").Value = -21223.4995
# Row 132
$ws.Range("H132").Value = 2694.4736
$ws.Range("I132").Value = 2258.5293
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 6775.5879
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -4245.5879
$ws.Range("N132").Value = -24260

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 64.666664
$ws.Range("J2").Value = 72
$ws.Range("L2").Value = 432
$ws.Range("N2").Value = -658
# Row 5
$ws.Range("H5").Value = 645.7826
$ws.Range("I5").Value = 415.33334
$ws.Range("K5").Value = 1246.00002
$ws.Range("M5").Value = -1134.00002
# Row 86
$ws.Range("H86").Value = 1606.1111
$ws.Range("I86").Value = 1612.8334
$ws.Range("J86").Value = 1592.6666
$ws.Range("K86").Value = 4838.5002
$ws.Range("L86").Value = 4777.9998
$ws.Range("M86").Value = -3652.5002
$ws.Range("N86").Value = -7149.9998
# Row 89
$ws.Range("H89").Value = 1606.1111
$ws.Range("I89").Value = 1612.8334
$ws.Range("J89").Value = 1592.6666
$ws.Range("K89").Value = 14515.5006
$ws.Range("L89").Value = 14333.9994
$ws.Range("M89").Value = -8587.500599999999
$ws.Range("N89").Value = -26189.9994
# Row 135
$ws.Range("H135").Value = 645.7826
$ws.Range("I135").Value = 415.33334
$ws.Range("K135").Value = 3738.00006
$ws.Range("M135").Value = -1203.00006
# Row 139
$ws.Range("H139").Value = 76927380
$ws.Range("I139").Value = 111113544
$ws.Range("J139").Value = 8500
$ws.Range("K139").Value = 333340632
$ws.Range("L139").Value = 25500
$ws.Range("M139").Value = -333335492
$ws.Range("N139").Value = -35780

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 2857269.2
$ws.Range("I2").Value = 39.285713
$ws.Range("J2").Value = 7143114
$ws.Range("K2").Value = 39.285713
$ws.Range("L2").Value = 7143114
$ws.Range("M2").Value = 73.714287
$ws.Range("N2").Value = -7143340
# Row 11
$ws.Range("H11").Value = 3531000
$ws.Range("I11").Value = 6500000
$ws.Range("K11").Value = 6500000
$ws.Range("M11").Value = -6499861
# Row 117
$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -51884
# Row 126
$ws.Range("H126").Value = 5381.5
$ws.Range("J126").Value = 5257
$ws.Range("L126").Value = 15771
$ws.Range("N126").Value = -20711
# Row 132
$ws.Range("H132").Value = 8145.9653
$ws.Range("I132").Value = 7249.32
$ws.Range("K132").Value = 21747.96
$ws.Range("M132").Value = -19217.96

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4504.0625
$ws.Range("I7").Value = 2703.2
$ws.Range("K7").Value = 2703.2
$ws.Range("M7").Value = -2591.2
# Row 61
$ws.Range("H61").Value = 1945.2858
$ws.Range("I61").Value = 2065
$ws.Range("K61").Value = 2065
$ws.Range("M61").Value = -1863
# Row 93
$ws.Range("H93").Value = 3695.077
$ws.Range("I93").Value = 3542.889
$ws.Range("J93").Value = 4037.5
$ws.Range("K93").Value = 3542.889
$ws.Range("L93").Value = 4037.5
$ws.Range("M93").Value = -2294.889
$ws.Range("N93").Value = -6533.5
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 113
$ws.Range("H113").Value = 1945.2858
$ws.Range("I113").Value = 2065
$ws.Range("K113").Value = 2065
$ws.Range("M113").Value = 105
# Row 122
$ws.Range("H122").Value = 4456.375
$ws.Range("I122").Value = 3650.2
$ws.Range("J122").Value = 5800
$ws.Range("K122").Value = 10950.6
$ws.Range("L122").Value = 17400
$ws.Range("M122").Value = -8500.599999999999
$ws.Range("N122").Value = -22300
# Row 126
$ws.Range("H126").Value = 4504.0625
$ws.Range("I126").Value = 2703.2
$ws.Range("K126").Value = 8109.599999999999
$ws.Range("M126").Value = -5639.599999999999
# Row 136
$ws.Range("H136").Value = 3814.2693
$ws.Range("I136").Value = 3298.1875
$ws.Range("J136").Value = 4640
$ws.Range("K136").Value = 9894.5625
$ws.Range("L136").Value = 13920
$ws.Range("M136").Value = -7344.5625
$ws.Range("N136").Value = -19020

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 790.1539
$ws.Range("I107").Value = 780.8182
$ws.Range("J107").Value = 841.5
$ws.Range("K107").Value = 2342.4546
$ws.Range("L107").Value = 2524.5
$ws.Range("M107").Value = -422.4546
$ws.Range("N107").Value = -6364.5
# Row 122
$ws.Range("H122").Value = 1439.5714
$ws.Range("I122").Value = 1431.1666
$ws.Range("K122").Value = 4293.4998
$ws.Range("M122").Value = -1843.4998
# Row 132
$ws.Range("H132").Value = 2610.9443
$ws.Range("I132").Value = 1771.2142
$ws.Range("K132").Value = 5313.642599999999
$ws.Range("M132").Value = -2783.642599999999
# Row 136
$ws.Range("H136").Value = 1362.3125
$ws.Range("I136").Value = 1328
$ws.Range("J136").Value = 1602.5
$ws.Range("K136").Value = 3984
$ws.Range("L136").Value = 4807.5
$ws.Range("M136").Value = -1434
$ws.Range("N136").Value = -9907.5
